# Apply updated symbol list values (Price / Volume(1h) columns) per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (no auto number/percent conversion),
# then restore the original (unstyled) cell formatting by copying the style
# from the never-touched "Hora" column (G) in the same row.
function Set-TextValue($ws, $addr, $rowNum, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = $ws.Cells.Item($rowNum, 7).Style
}

Set-TextValue $ws "D2" 2 "309.73"
Set-TextValue $ws "E2" 2 "-3.71%"
Set-TextValue $ws "D3" 3 "48.53"
Set-TextValue $ws "E3" 3 "-5.33%"
Set-TextValue $ws "D4" 4 "5.166"
Set-TextValue $ws "E4" 4 "-3.45%"
Set-TextValue $ws "D5" 5 "0.07768"
Set-TextValue $ws "E5" 5 "-4.21%"
Set-TextValue $ws "D6" 6 "4.482"
Set-TextValue $ws "E6" 6 "-2.11%"
Set-TextValue $ws "D7" 7 "1.317"
Set-TextValue $ws "E7" 7 "16.84%"
Set-TextValue $ws "D8" 8 "1.560"
Set-TextValue $ws "E8" 8 "-5.60%"
Set-TextValue $ws "D9" 9 "0.1222"
Set-TextValue $ws "E9" 9 "-7.20%"
Set-TextValue $ws "D10" 10 "0.1948"
Set-TextValue $ws "E10" 10 "-0.34%"
Set-TextValue $ws "D11" 11 "0.04684"
Set-TextValue $ws "E11" 11 "2.40%"
Set-TextValue $ws "D12" 12 "0.09392"
Set-TextValue $ws "E12" 12 "-2.66%"
Set-TextValue $ws "E13" 13 "0.16%"
Set-TextValue $ws "D14" 14 "0.001272"
Set-TextValue $ws "E14" 14 "-3.44%"
Set-TextValue $ws "D15" 15 "0.04158"
Set-TextValue $ws "E15" 15 "-3.46%"
Set-TextValue $ws "D16" 16 "0.005814"
Set-TextValue $ws "E16" 16 "0.16%"
Set-TextValue $ws "D17" 17 "3.330"
Set-TextValue $ws "E17" 17 "-1.46%"
Set-TextValue $ws "D18" 18 "2.274"
Set-TextValue $ws "E18" 18 "-6.50%"
Set-TextValue $ws "E19" 19 "2.93%"
Set-TextValue $ws "D20" 20 "8.310"
Set-TextValue $ws "E20" 20 "1.41%"
Set-TextValue $ws "D21" 21 "0.1347"
Set-TextValue $ws "E21" 21 "-3.59%"
Set-TextValue $ws "E22" 22 "3.60%"
Set-TextValue $ws "D23" 23 "0.001278"
Set-TextValue $ws "E23" 23 "-2.12%"
Set-TextValue $ws "D24" 24 "0.004106"
Set-TextValue $ws "E24" 24 "-4.70%"
Set-TextValue $ws "E25" 25 "0.04%"
Set-TextValue $ws "E26" 26 "-3.88%"
Set-TextValue $ws "E38" 38 "-7.12%"
Set-TextValue $ws "D39" 39 "0.05898"
Set-TextValue $ws "E39" 39 "6.39%"
Set-TextValue $ws "E40" 40 "70.97%"
Set-TextValue $ws "D41" 41 "0.007927"
Set-TextValue $ws "E41" 41 "1.97%"
Set-TextValue $ws "D42" 42 "0.1424"
Set-TextValue $ws "E42" 42 "-1.38%"
Set-TextValue $ws "D43" 43 "0.008416"
Set-TextValue $ws "E43" 43 "9.64%"
Set-TextValue $ws "D44" 44 "0.007667"
Set-TextValue $ws "E44" 44 "-13.15%"
Set-TextValue $ws "D45" 45 "0.3101"
Set-TextValue $ws "E45" 45 "-12.03%"
Set-TextValue $ws "D46" 46 "0.00006895"
Set-TextValue $ws "E46" 46 "1.09%"
Set-TextValue $ws "E47" 47 "0.16%"
Set-TextValue $ws "D48" 48 "0.05670"
Set-TextValue $ws "E48" 48 "-7.42%"
Set-TextValue $ws "E49" 49 "0.32%"
Set-TextValue $ws "E50" 50 "0.16%"
Set-TextValue $ws "D51" 51 "0.0002000"
Set-TextValue $ws "E51" 51 "0.16%"
